$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New gear rows appended below the existing table (rows 11-13)
$newRows = @(
    @{ Row = 11; A = 2406038; B = "C 0624129"; C = "2E24"; D = "2E24"; E = "2D24"; F = "2E24" },
    @{ Row = 12; A = 2406039; B = "C 0624129"; C = "2E24"; D = "2E24"; E = "2D24"; F = "2E24" },
    @{ Row = 13; A = 2406040; B = "C 0624129"; C = "2E24"; D = "2E24"; E = "2D24"; F = "2E24" }
)

foreach ($r in $newRows) {
    # Give the new row the same look (borders/number formats/fonts) as the row above it
    $ws.Range("A10:F10").Copy() | Out-Null
    $ws.Range("A$($r.Row):F$($r.Row)").PasteSpecial(-4122) | Out-Null
    $ws.Rows.Item($r.Row).RowHeight = $ws.Rows.Item(10).RowHeight

    $ws.Cells.Item($r.Row, 1).Value = $r.A
    $ws.Cells.Item($r.Row, 2).Value = $r.B
    $ws.Cells.Item($r.Row, 3).Value = $r.C
    $ws.Cells.Item($r.Row, 4).Value = $r.D
    $ws.Cells.Item($r.Row, 5).Value = $r.E
    $ws.Cells.Item($r.Row, 6).Value = $r.F
}

$excel.CutCopyMode = $false

# Move the selection to match the final state (single cell A13 selected)
$ws.Range("A13").Select()
